$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing "Date" values in column B to reflect the new test run timestamps
$ws.Range("B2").Value = "Thu Jan 25 17:47:55 EST 2024"
$ws.Range("B4").Value = "Thu Jan 25 17:48:10 EST 2024"

# Add new row for Estate Tax test data
$ws.Range("D5").Value = "New Tax Return Amount Due"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("E5").Value = "Estate Tax"

# Move selection to the newly added cell, matching the saved workbook state
$ws.Range("E5").Select()
